# Update the cryptocurrency price/volume snapshot (GitHub Actions refresh).
#
# Only the Price (column D) and Volume(1h) (column E) cells change value,
# except for rows 14/15 where ShibaInu and Polygon swap rank positions
# (their Coin name, Link, Price and Volume cells all change).
#
# Price cells are plain text (e.g. thousands are dot-separated like
# "29.436.29", which is not a valid number). Whenever the new price text
# *would* parse as an ordinary decimal number (e.g. "240.71"), the cell's
# NumberFormat is set to Text ("@") first so the COM layer keeps storing it
# as a string instead of silently converting it to a numeric value - this
# mirrors how the original values were authored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '29.436.29'
$ws.Range("E2").Value = '  +0.37%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.848.76'
$ws.Range("E3").Value = '  +0.37%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.12%  '

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.71'
$ws.Range("E5").Value = '  +0.84%  '

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6299'
$ws.Range("E6").Value = '  -0.06%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  +0.08%  '

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07676'
$ws.Range("E8").Value = '  +1.96%  '

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2924'

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.68'
$ws.Range("E10").Value = '  +0.72%  '

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07744'
$ws.Range("E11").Value = '  +0.62%  '

# Row 12 - WrappedEther
$ws.Range("D12").Value = '1.861.85'
$ws.Range("E12").Value = '  +1.09%  '

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.030'
$ws.Range("E13").Value = '  +1.02%  '

# Row 14 - now Polygon (was ShibaInu)
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6791'
$ws.Range("E14").Value = '  +0.01%  '

# Row 15 - now ShibaInu (was Polygon)
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001071'
$ws.Range("E15").Value = '  +2.39%  '

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.72'
$ws.Range("E16").Value = '  +0.86%  '

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = '2.115.61'
$ws.Range("E17").Value = '  +1.61%  '

# Row 18 - Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.194'
$ws.Range("E18").Value = '  +0.80%  '

# Row 19
$ws.Range("D19").Value = '29.466.99'
$ws.Range("E19").Value = '  +0.40%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.45'
$ws.Range("E20").Value = '  -0.06%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.43'
$ws.Range("E21").Value = '  +0.07%  '

# Row 22
$ws.Range("E22").Value = '  +0.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.434'
$ws.Range("E23").Value = '  +0.19%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.85'
$ws.Range("E25").Value = '  +0.91%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1378'
$ws.Range("E26").Value = '  -0.97%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.413'

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.69'
$ws.Range("E28").Value = '  +0.53%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.344'
$ws.Range("E29").Value = '  +5.97%  '

# Row 30
$ws.Range("E30").Value = '  +0.60%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05665'
$ws.Range("E31").Value = '  +0.53%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.120'
$ws.Range("E32").Value = '  +0.36%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.026'
$ws.Range("E33").Value = '  +0.16%  '

# Row 34
$ws.Range("E34").Value = '  +0.73%  '

# Row 35 - ARBITRUM
$ws.Range("E35").Value = '  +0.60%  '

# Row 36 - ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7050'
$ws.Range("E36").Value = '  -0.60%  '

# Row 37 - HuobiToken
$ws.Range("E37").Value = '  -0.22%  '

# Row 38 - MXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.780'
$ws.Range("E38").Value = '  +0.90%  '

# Row 39 - VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01792'
$ws.Range("E39").Value = '  -1.00%  '

# Row 40 - Maker
$ws.Range("D40").Value = '1.219.60'
$ws.Range("E40").Value = '  -1.98%  '

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.539'
$ws.Range("E41").Value = '  +4.82%  '

# Row 42 - TrustWalletToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9102'
$ws.Range("E42").Value = '  +0.96%  '

# Row 43 - PaxDollar
$ws.Range("E43").Value = '  +0.19%  '

# Row 44 - Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.65'
$ws.Range("E44").Value = '  -0.21%  '

# Row 45 - Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.11'
$ws.Range("E45").Value = '  +0.69%  '

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = '  +1.91%  '

# Row 47 - Aptos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.129'
$ws.Range("E47").Value = '  +0.14%  '

# Row 48 - TheSandbox
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4022'
$ws.Range("E48").Value = '  +0.74%  '

# Row 49 - EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.998'
$ws.Range("E49").Value = '  +0.60%  '

# Row 50 - RenderToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.679'
$ws.Range("E50").Value = '  +0.44%  '

# Row 51 - Algorand
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1144'
$ws.Range("E51").Value = '  +2.10%  '
